# Auto-generated Excel COM-interop script
# Reverts 'Product' placeholder text back from the multi-industry 'AI/ML' template
# (restores AI/ML Training Schedule content -> Product-branded content)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Training Schedule Overview ---
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Product Development IMPLEMENTATION PROJECT - TRProductNING SCHEDULE"
$ws.Range("C4").Value = "Enterprise Product Development Implementation"
$ws.Range("A7").Value = "TRProductNING SCHEDULE SUMMARY"
$ws.Range("A9").Value = "Product Development Fundamentals (Product-101)"
$ws.Range("A10").Value = "Product Development Platform Overview (Product-102)"
$ws.Range("A11").Value = "Data Analysis for Business (Product-201)"
$ws.Range("A12").Value = "Advanced Product Techniques (Product-301)"
$ws.Range("A13").Value = "ProductOps for IT Teams (Product-302)"
$ws.Range("B13").Value = "Product Engineers, IT"
$ws.Range("A14").Value = "Model Validation & QA (Product-303)"
$ws.Range("B14").Value = "Product Engineers, QA"
$ws.Range("A15").Value = "Executive Overview (Product-401)"
$ws.Range("A16").Value = "Train-the-Trainer (Product-501)"
$ws.Range("A18").Value = "TRProductNING SCHEDULE STATISTICS"

# Insert blank (empty) rows that exist in the row index but carry no cell data
$ws.Rows.Item(3).OutlineLevel = 0
$ws.Rows.Item(6).OutlineLevel = 0
$ws.Rows.Item(17).OutlineLevel = 0

# --- Sheet 2: Detailed Training Schedule ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("A1").Value = "DETProductLED TRProductNING SCHEDULE"
$ws.Range("A4").Value = "Product-101"
$ws.Range("B4").Value = "Product Development Fundamentals"
$ws.Range("A5").Value = "Product-102"
$ws.Range("B5").Value = "Product Development Platform Overview"
$ws.Range("A6").Value = "Product-201"
$ws.Range("A7").Value = "Product-201"
$ws.Range("A8").Value = "Product-201"
$ws.Range("A9").Value = "Product-301"
$ws.Range("B9").Value = "Advanced Product Techniques"
$ws.Range("A10").Value = "Product-301"
$ws.Range("B10").Value = "Advanced Product Techniques"
$ws.Range("A11").Value = "Product-302"
$ws.Range("B11").Value = "ProductOps for IT Teams"
$ws.Range("C11").Value = "Product Engineers, IT"
$ws.Range("A12").Value = "Product-302"
$ws.Range("B12").Value = "ProductOps for IT Teams"
$ws.Range("C12").Value = "Product Engineers, IT"
$ws.Range("A13").Value = "Product-303"
$ws.Range("C13").Value = "Product Engineers, QA"
$ws.Range("A14").Value = "Product-303"
$ws.Range("C14").Value = "Product Engineers, QA"
$ws.Range("A15").Value = "Product-401"
$ws.Range("A16").Value = "Product-501"
$ws.Range("A17").Value = "Product-501"
$ws.Range("A18").Value = "Product-501"
$ws.Range("A19").Value = "Product-501"
$ws.Range("A20").Value = "Product-501"

# Insert blank (empty) rows that exist in the row index but carry no cell data
$ws.Rows.Item(2).OutlineLevel = 0

# --- Sheet 3: Instructor Schedule ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("B7").Value = "Advanced Product Techniques"
$ws.Range("B8").Value = "Advanced Product Techniques"
$ws.Range("B9").Value = "ProductOps for IT Teams"
$ws.Range("B10").Value = "ProductOps for IT Teams"

# Insert blank (empty) rows that exist in the row index but carry no cell data
$ws.Rows.Item(2).OutlineLevel = 0

# --- Sheet 4: Facility Schedule ---
$ws = $wb.Worksheets.Item(4)

$ws.Range("B7").Value = "Advanced Product Techniques"
$ws.Range("B8").Value = "Advanced Product Techniques"
$ws.Range("B9").Value = "ProductOps for IT Teams"
$ws.Range("B10").Value = "ProductOps for IT Teams"

# Insert blank (empty) rows that exist in the row index but carry no cell data
$ws.Rows.Item(2).OutlineLevel = 0

# --- Sheet 5: Participant Tracking ---
$ws = $wb.Worksheets.Item(5)

$ws.Range("F4").Value = "Product-101"
$ws.Range("F5").Value = "Product-102"
$ws.Range("F6").Value = "Product-301"
$ws.Range("F7").Value = "Product-302"
$ws.Range("E8").Value = "Product Engineer"
$ws.Range("F8").Value = "Product-101"
$ws.Range("E9").Value = "Product Engineer"
$ws.Range("F9").Value = "Product-102"
$ws.Range("E10").Value = "Product Engineer"
$ws.Range("F10").Value = "Product-302"
$ws.Range("E11").Value = "Product Engineer"
$ws.Range("F11").Value = "Product-303"
$ws.Range("F12").Value = "Product-101"
$ws.Range("F13").Value = "Product-102"
$ws.Range("F14").Value = "Product-401"
$ws.Range("F15").Value = "Product-101"
$ws.Range("F16").Value = "Product-102"
$ws.Range("F17").Value = "Product-501"
$ws.Range("F18").Value = "Product-101"
$ws.Range("F19").Value = "Product-102"
$ws.Range("F20").Value = "Product-301"
$ws.Range("F21").Value = "Product-303"
$ws.Range("F22").Value = "Product-501"

# Insert blank (empty) rows that exist in the row index but carry no cell data
$ws.Rows.Item(2).OutlineLevel = 0

